$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force only the numeric-looking data columns to Text format so values like
# "225.00" or "13" are preserved verbatim as text instead of becoming numbers
# (the other columns are never numeric-looking, so Excel keeps them text by default).
$ws.Range("G2:K13").NumberFormat = "@"

# Header row (row 1), columns A-K
$ws.Range("A1").Value = "venue"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "result"
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("F1").Value = "batsman"
$ws.Range("G1").Value = "totalRuns"
$ws.Range("H1").Value = "totalBalls"
$ws.Range("I1").Value = "total4s"
$ws.Range("J1").Value = "total6s"
$ws.Range("K1").Value = "sr"

# Row 2
$ws.Range("A2").Value = " Abu Dhabi"
$ws.Range("B2").Value = " October 25 2020"
$ws.Range("C2").Value = "Royals won by 8 wickets (with 10 balls remaining)"
$ws.Range("D2").Value = "Rajasthan Royals"
$ws.Range("E2").Value = "Mumbai Indians"
$ws.Range("F2").Value = "Robin Uthappa "
$ws.Range("G2").Value = "13"
$ws.Range("H2").Value = "11"
$ws.Range("I2").Value = "2"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "118.18"

# Row 3
$ws.Range("A3").Value = " Dubai (DSC)"
$ws.Range("B3").Value = " October 17 2020"
$ws.Range("C3").Value = "RCB won by 7 wickets (with 2 balls remaining)"
$ws.Range("D3").Value = "Rajasthan Royals"
$ws.Range("E3").Value = "Royal Challengers Bangalore"
$ws.Range("F3").Value = "Robin Uthappa "
$ws.Range("G3").Value = "41"
$ws.Range("H3").Value = "22"
$ws.Range("I3").Value = "7"
$ws.Range("J3").Value = "1"
$ws.Range("K3").Value = "186.36"

# Row 4
$ws.Range("A4").Value = " Sharjah"
$ws.Range("B4").Value = " September 27 2020"
$ws.Range("C4").Value = "Royals won by 4 wickets (with 3 balls remaining)"
$ws.Range("D4").Value = "Rajasthan Royals"
$ws.Range("E4").Value = "Kings XI Punjab"
$ws.Range("F4").Value = "Robin Uthappa "
$ws.Range("G4").Value = "9"
$ws.Range("H4").Value = "4"
$ws.Range("I4").Value = "2"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "225.00"

# Row 5
$ws.Range("A5").Value = " Abu Dhabi"
$ws.Range("B5").Value = " October 30 2020"
$ws.Range("C5").Value = "Royals won by 7 wickets (with 15 balls remaining)"
$ws.Range("D5").Value = "Rajasthan Royals"
$ws.Range("E5").Value = "Kings XI Punjab"
$ws.Range("F5").Value = "Robin Uthappa "
$ws.Range("G5").Value = "30"
$ws.Range("H5").Value = "23"
$ws.Range("I5").Value = "1"
$ws.Range("J5").Value = "2"
$ws.Range("K5").Value = "130.43"

# Row 6
$ws.Range("A6").Value = " Dubai (DSC)"
$ws.Range("B6").Value = " September 30 2020"
$ws.Range("C6").Value = "KKR won by 37 runs"
$ws.Range("D6").Value = "Rajasthan Royals"
$ws.Range("E6").Value = "Kolkata Knight Riders"
$ws.Range("F6").Value = "Robin Uthappa "
$ws.Range("G6").Value = "2"
$ws.Range("H6").Value = "7"
$ws.Range("I6").Value = "0"
$ws.Range("J6").Value = "0"
$ws.Range("K6").Value = "28.57"

# Row 7
$ws.Range("A7").Value = " Abu Dhabi"
$ws.Range("B7").Value = " October 03 2020"
$ws.Range("C7").Value = "RCB won by 8 wickets (with 5 balls remaining)"
$ws.Range("D7").Value = "Rajasthan Royals"
$ws.Range("E7").Value = "Royal Challengers Bangalore"
$ws.Range("F7").Value = "Robin Uthappa "
$ws.Range("G7").Value = "17"
$ws.Range("H7").Value = "22"
$ws.Range("I7").Value = "1"
$ws.Range("J7").Value = "0"
$ws.Range("K7").Value = "77.27"

# Row 8
$ws.Range("A8").Value = " Dubai (DSC)"
$ws.Range("B8").Value = " October 14 2020"
$ws.Range("C8").Value = "Capitals won by 13 runs"
$ws.Range("D8").Value = "Rajasthan Royals"
$ws.Range("E8").Value = "Delhi Capitals"
$ws.Range("F8").Value = "Robin Uthappa "
$ws.Range("G8").Value = "32"
$ws.Range("H8").Value = "27"
$ws.Range("I8").Value = "3"
$ws.Range("J8").Value = "1"
$ws.Range("K8").Value = "118.51"

# Row 9
$ws.Range("A9").Value = " Dubai (DSC)"
$ws.Range("B9").Value = " October 22 2020"
$ws.Range("C9").Value = "Sunrisers won by 8 wickets (with 11 balls remaining)"
$ws.Range("D9").Value = "Rajasthan Royals"
$ws.Range("E9").Value = "Sunrisers Hyderabad"
$ws.Range("F9").Value = "Robin Uthappa "
$ws.Range("G9").Value = "19"
$ws.Range("H9").Value = "13"
$ws.Range("I9").Value = "2"
$ws.Range("J9").Value = "1"
$ws.Range("K9").Value = "146.15"

# Row 10
$ws.Range("A10").Value = " Dubai (DSC)"
$ws.Range("B10").Value = " November 01 2020"
$ws.Range("C10").Value = "KKR won by 60 runs"
$ws.Range("D10").Value = "Rajasthan Royals"
$ws.Range("E10").Value = "Kolkata Knight Riders"
$ws.Range("F10").Value = "Robin Uthappa "
$ws.Range("G10").Value = "6"
$ws.Range("H10").Value = "2"
$ws.Range("I10").Value = "0"
$ws.Range("J10").Value = "1"
$ws.Range("K10").Value = "300.00"

# Row 11
$ws.Range("A11").Value = " Sharjah"
$ws.Range("B11").Value = " September 22 2020"
$ws.Range("C11").Value = "Royals won by 16 runs"
$ws.Range("D11").Value = "Rajasthan Royals"
$ws.Range("E11").Value = "Chennai Super Kings"
$ws.Range("F11").Value = "Robin Uthappa "
$ws.Range("G11").Value = "5"
$ws.Range("H11").Value = "9"
$ws.Range("I11").Value = "0"
$ws.Range("J11").Value = "0"
$ws.Range("K11").Value = "55.55"

# Row 12
$ws.Range("A12").Value = " Abu Dhabi"
$ws.Range("B12").Value = " October 19 2020"
$ws.Range("C12").Value = "Royals won by 7 wickets (with 15 balls remaining)"
$ws.Range("D12").Value = "Rajasthan Royals"
$ws.Range("E12").Value = "Chennai Super Kings"
$ws.Range("F12").Value = "Robin Uthappa "
$ws.Range("G12").Value = "4"
$ws.Range("H12").Value = "9"
$ws.Range("I12").Value = "0"
$ws.Range("J12").Value = "0"
$ws.Range("K12").Value = "44.44"

# Row 13
$ws.Range("A13").Value = " Dubai (DSC)"
$ws.Range("B13").Value = " October 11 2020"
$ws.Range("C13").Value = "Royals won by 5 wickets (with 1 ball remaining)"
$ws.Range("D13").Value = "Rajasthan Royals"
$ws.Range("E13").Value = "Sunrisers Hyderabad"
$ws.Range("F13").Value = "Robin Uthappa "
$ws.Range("G13").Value = "18"
$ws.Range("H13").Value = "15"
$ws.Range("I13").Value = "1"
$ws.Range("J13").Value = "1"
$ws.Range("K13").Value = "120.00"
